$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 2008年 and 2009年 rows (rows 2 and 3).
# This shifts 2010年 (was row 4) up to row 2, and 2011年 (was row 5) up to row 3.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
